$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '36.902.06'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  -0.60%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '2.038.15'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  -0.71%  '
$ws.Range("E4").Value = '  +0.04%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '244.97'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -1.71%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.655'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -1.28%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '58.16'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  -0.71%  '
$ws.Range("E8").Value = '  -0.05%  '
$ws.Range("E9").Value = '  -1.77%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.0768'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  -1.89%  '
$ws.Range("E11").Value = '  +2.19%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '15.30'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  -3.42%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.878'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  +7.88%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '2.335.94'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  -0.69%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '5.60'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  +0.78%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '2.041.27'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  -0.61%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '18.10'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  +7.96%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '36.876.65'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  -0.76%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '73.48'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  -1.61%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '0.0₃0884'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  -1.99%  '
$ws.Range("E21").Value = '  -0.05%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '234.83'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -0.70%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '2.44'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +2.34%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '9.54'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +4.01%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '168.10'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  -0.03%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '2.12'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -3.04%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '19.87'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +0.21%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '5.43'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  +15.25%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '0.123'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  -0.44%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '1.10'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -3.33%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '4.68'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +4.21%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.0610'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  -0.48%  '
$ws.Range("E34").Value = '  +0.00%  '
$ws.Range("E35").Value = '  -3.73%  '
$ws.Range("E36").Value = '  +6.05%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '2.23'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  +0.36%  '
$ws.Range("E38").Value = '  -4.11%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '3.12'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -1.69%  '
$ws.Range("E40").Value = '  -0.37%  '
$ws.Range("E41").Value = '  +0.19%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.0956'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -11.69%  '
$ws.Range("E43").Value = '  +0.87%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '96.86'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +0.97%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '16.78'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -4.61%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '1.289.78'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +0.68%  '
$ws.Range("E47").Value = '  -4.36%  '
$ws.Range("E48").Value = '  -0.41%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '3.70'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +7.64%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '6.68'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -1.08%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '2.222.94'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -0.81%  '
